$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the floating-point representation of the timestamp in A16
$ws.Range("A16").Value = 45868.70855333333

# Append the new row of sensor data (row 17)
$ws.Range("A17").Value = 45868.75023297059
$ws.Range("B17").Value = 2025
$ws.Range("C17").Value = 31
$ws.Range("D17").Value = 15.56
$ws.Range("E17").Value = 87.05
$ws.Range("F17").Value = 8.789999999999999
$ws.Range("G17").Value = 11.98
$ws.Range("H17").Value = "ESE"
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = "18:00:20"

# Match the date number format used by the other timestamp cells in column A
$ws.Range("A17").NumberFormat = $ws.Range("A16").NumberFormat
